$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the "QnA Number" column (column A) ---
# Shift columns B:C left into A:B (values, then formats) using copy/paste,
# which correctly carries over both content and cell styles.
$ws.Range("B2:C17").Copy()
$ws.Range("A2").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("B2:C17").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Clear out the now-stale column C (its data now lives in column B)
$ws.Range("C2:C17").Clear()

# --- Step 2: remove the now-empty row 1 ---
# Shift rows 2:17 up into rows 1:16 (values, then formats).
$ws.Range("A2:B17").Copy()
$ws.Range("A1").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("A2:B17").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Clear the now-stale last row (row 17)
$ws.Range("A17:B17").Clear()

# --- Column widths for the new A/B layout ---
$ws.Columns("A").ColumnWidth = 45.28515625
$ws.Columns("B").ColumnWidth = 100.85546875

# Reset the selection to A1 (default)
$ws.Range("A1").Select()
